$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2260.7693
$ws.Range("I32").Value = 1866.3334
$ws.Range("J32").Value = 2379.1
$ws.Range("K32").Value = 1866.3334
$ws.Range("L32").Value = 2379.1
$ws.Range("M32").Value = -1540.3334
$ws.Range("N32").Value = -3031.1

$ws.Range("H38").Value = 1591.5294
$ws.Range("I38").Value = 858.38464
$ws.Range("J38").Value = 3974.25
$ws.Range("K38").Value = 2575.15392
$ws.Range("L38").Value = 11922.75
$ws.Range("M38").Value = -2203.15392
$ws.Range("N38").Value = -12666.75

$ws.Range("H39").Value = 126
$ws.Range("I39").Value = 82.5
$ws.Range("K39").Value = 247.5
$ws.Range("M39").Value = 48.5

$ws.Range("H40").Value = 3385.3333
$ws.Range("I40").Value = 2141.5
$ws.Range("J40").Value = 4380.4
$ws.Range("K40").Value = 2141.5
$ws.Range("L40").Value = 4380.4
$ws.Range("M40").Value = -1966.5
$ws.Range("N40").Value = -4730.4

$ws.Range("H42").Value = 52631884
$ws.Range("I42").Value = 66666890
$ws.Range("J42").Value = 612
$ws.Range("K42").Value = 200000670
$ws.Range("L42").Value = 1836
$ws.Range("M42").Value = -200000440
$ws.Range("N42").Value = -2296

$ws.Range("H43").Value = 5628.3213
$ws.Range("J43").Value = 6625.75
$ws.Range("L43").Value = 6625.75
$ws.Range("N43").Value = -6763.75

$ws.Range("H98").Value = 23713.857
$ws.Range("J98").Value = 999
$ws.Range("L98").Value = 999
$ws.Range("N98").Value = -3995

$ws.Range("H113").Value = 4425.25
$ws.Range("I113").Value = 2566
$ws.Range("J113").Value = 10003
$ws.Range("K113").Value = 2566
$ws.Range("L113").Value = 10003
$ws.Range("M113").Value = 688
$ws.Range("N113").Value = -16511

$ws.Range("H122").Value = 23713.857
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897

$ws.Range("H136").Value = 110000
$ws.Range("J136").Value = 110000
$ws.Range("L136").Value = 110000
$ws.Range("N136").Value = -120200

$ws.Range("H137").Value = 1427.3043
$ws.Range("I137").Value = 1239.0667
$ws.Range("K137").Value = 3717.2001
$ws.Range("M137").Value = -1167.2001

$ws.Range("H141").Value = 3779.375
$ws.Range("I141").Value = 3459.3076
$ws.Range("K141").Value = 10377.9228
$ws.Range("M141").Value = -5197.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 3500
$ws.Range("J14").Value = 3500
$ws.Range("L14").Value = 3500
$ws.Range("N14").Value = -3850

$ws.Range("H32").Value = 2763.9592
$ws.Range("I32").Value = 2488.889
$ws.Range("K32").Value = 2488.889
$ws.Range("M32").Value = -2201.889

$ws.Range("H34").Value = 119165.914
$ws.Range("I34").Value = 200000
$ws.Range("J34").Value = 111817.37
$ws.Range("K34").Value = 200000
$ws.Range("L34").Value = 111817.37
$ws.Range("M34").Value = -199729
$ws.Range("N34").Value = -112359.37

$ws.Range("H61").Value = 3922.5
$ws.Range("I61").Value = 4106.364
$ws.Range("K61").Value = 4106.364
$ws.Range("M61").Value = -3894.364

$ws.Range("H136").Value = 3922.5
$ws.Range("I136").Value = 4106.364
$ws.Range("K136").Value = 12319.092
$ws.Range("M136").Value = -9769.091999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2111.875
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 2127.8572
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2127.8572
$ws.Range("M20").Value = -1753
$ws.Range("N20").Value = -2621.8572

$ws.Range("H94").Value = 2014.6923
$ws.Range("I94").Value = 1449
$ws.Range("K94").Value = 1449
$ws.Range("M94").Value = -998

$ws.Range("H105").Value = 2872.8462
$ws.Range("I105").Value = 2872.8462
$ws.Range("K105").Value = 2872.8462
$ws.Range("M105").Value = -1125.8462

$ws.Range("H134").Value = 2229.8262
$ws.Range("I134").Value = 2214.3
$ws.Range("J134").Value = 2333.3333
$ws.Range("K134").Value = 6642.900000000001
$ws.Range("L134").Value = 6999.999899999999
$ws.Range("M134").Value = -4107.900000000001
$ws.Range("N134").Value = -12069.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 999
$ws.Range("I16").Value = 999
$ws.Range("K16").Value = 999
$ws.Range("M16").Value = -712

$ws.Range("H80").Value = 34974
$ws.Range("J80").Value = 34974
$ws.Range("L80").Value = 34974
$ws.Range("N80").Value = -37220

$ws.Range("H83").Value = 34974
$ws.Range("J83").Value = 34974
$ws.Range("L83").Value = 104922
$ws.Range("N83").Value = -116154

$ws.Range("H86").Value = 4335.3335
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4503
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4503
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6749

$ws.Range("H89").Value = 4335.3335
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4503
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 22515
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -33747

$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 999
$ws.Range("K113").Value = 999
$ws.Range("M113").Value = 1171

$ws.Range("H132").Value = 2478.923
$ws.Range("I132").Value = 2560.5
$ws.Range("K132").Value = 7681.5
$ws.Range("M132").Value = -5151.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 963.0909
$ws.Range("J38").Value = 1749.1666
$ws.Range("L38").Value = 5247.4998
$ws.Range("N38").Value = -5941.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5883.647
$ws.Range("I80").Value = 4688.143
$ws.Range("J80").Value = 6720.5
$ws.Range("K80").Value = 4688.143
$ws.Range("L80").Value = 6720.5
$ws.Range("M80").Value = -3690.143
$ws.Range("N80").Value = -8716.5

$ws.Range("H83").Value = 5883.647
$ws.Range("I83").Value = 4688.143
$ws.Range("J83").Value = 6720.5
$ws.Range("K83").Value = 23440.715
$ws.Range("L83").Value = 33602.5
$ws.Range("M83").Value = -18448.715
$ws.Range("N83").Value = -43586.5

$ws.Range("H136").Value = 35480
$ws.Range("J136").Value = 35480
$ws.Range("L136").Value = 106440
$ws.Range("N136").Value = -111540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4678
$ws.Range("I10").Value = 600
$ws.Range("J10").Value = 7736.5
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 7736.5
$ws.Range("M10").Value = -460
$ws.Range("N10").Value = -8016.5

$ws.Range("H40").Value = 6606.25
$ws.Range("I40").Value = 5240.385
$ws.Range("K40").Value = 5240.385
$ws.Range("M40").Value = -5104.385

$ws.Range("H93").Value = 3720.8
$ws.Range("J93").Value = 5394.316
$ws.Range("L93").Value = 5394.316
$ws.Range("N93").Value = -7890.316

$ws.Range("H96").Value = 35332.832
$ws.Range("J96").Value = 35332.832
$ws.Range("L96").Value = 35332.832
$ws.Range("N96").Value = -40824.832

$ws.Range("H100").Value = 6349.8335
$ws.Range("J100").Value = 7142.7856
$ws.Range("L100").Value = 7142.7856
$ws.Range("N100").Value = -8224.785599999999

$ws.Range("H122").Value = 7114.5293
$ws.Range("I122").Value = 7451.4287
$ws.Range("J122").Value = 5542.3335
$ws.Range("K122").Value = 22354.2861
$ws.Range("L122").Value = 16627.0005
$ws.Range("M122").Value = -19904.2861
$ws.Range("N122").Value = -21527.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H81").Value = 4998.364
$ws.Range("I81").Value = 4133.3335
$ws.Range("J81").Value = 6036.4
$ws.Range("K81").Value = 8266.666999999999
$ws.Range("L81").Value = 12072.8
$ws.Range("M81").Value = -7205.666999999999
$ws.Range("N81").Value = -14194.8

$ws.Range("H84").Value = 4998.364
$ws.Range("I84").Value = 4133.3335
$ws.Range("J84").Value = 6036.4
$ws.Range("K84").Value = 41333.335
$ws.Range("L84").Value = 60364
$ws.Range("M84").Value = -36029.335
$ws.Range("N84").Value = -70972

$ws.Range("H100").Value = 1765.091
$ws.Range("I100").Value = 1774.5555
$ws.Range("K100").Value = 3549.111
$ws.Range("M100").Value = -3008.111

$ws.Range("H119").Value = 68332.664
$ws.Range("J119").Value = 68332.664
$ws.Range("L119").Value = 68332.664
$ws.Range("N119").Value = -78008.664

$ws.Range("H126").Value = 1143.7142
$ws.Range("I126").Value = 1182.1666
$ws.Range("K126").Value = 3546.4998
$ws.Range("M126").Value = -1076.4998

$ws.Range("H137").Value = 56000
$ws.Range("I137").Value = 56000
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 56000
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -50900
$ws.Range("N137").ClearContents()
